$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.456.93"
$ws.Range("E2").Value = "  -3.09%  "
$ws.Range("D3").Value = "1.669.39"
$ws.Range("E3").Value = "  -2.19%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.86"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5156"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06471"
$ws.Range("E8").Value = "  -1.65%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2578"
$ws.Range("E9").Value = "  -2.83%  "
$ws.Range("E10").Value = "  -3.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07680"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.350"
$ws.Range("E12").Value = "  -4.72%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.674.95"
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "1.898.32"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5565"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("D16").Value = "0.0₅8063"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.81"
$ws.Range("E17").Value = "  -4.26%  "
$ws.Range("D18").Value = "26.475.29"
$ws.Range("E18").Value = "  -3.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.72"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.430"
$ws.Range("E21").Value = "  -5.00%  "
$ws.Range("E22").Value = "  -2.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.899"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.17"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.738"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1164"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.018"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05221"
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.261"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.372"
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.222"
$ws.Range("E33").Value = "  -5.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.587"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.758"
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9248"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5737"
$ws.Range("E38").Value = "  -2.00%  "
$ws.Range("D39").Value = "1.166.41"
$ws.Range("E39").Value = "  +11.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01600"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.007"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8419"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.647"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.24"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "1.808.26"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -3.07%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4501"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.14"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.937"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05119"
$ws.Range("E51").Value = "  -2.33%  "
